$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.261.95"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.72%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.906.51"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.80%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.50"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5413"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +3.64%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3816"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07296"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.22"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +5.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9043"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.56%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08186"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "96.03"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.357"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.9995"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.90"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008653"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9999"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "27.288.70"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.208.81"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -37.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.055"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.81"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.524"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.90%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "148.59"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.311"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.39"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.40%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.752"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "116.99"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.868"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.677"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09217"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8319"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +5.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05082"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.225"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.012"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.323"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.75%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.704"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +4.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6001"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +4.77%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02005"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.082"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.292"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.671"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.89%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "116.43"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5158"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +6.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1532"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.22"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9991"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.647"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.70%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "38.29"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.56%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.89%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.70"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.35%  "
